$wb = $excel.ActiveWorkbook

# 1. Add defined name CO2_to_C -> conversion!$B$5
$wb.Names.Add("CO2_to_C", "=conversion!`$B`$5")

$ws = $wb.Worksheets.Item("Emission_Coefficient")

# --- Header row (order chosen to match shared-string introduction order) ---
$ws.Range("H1").Value = "EC (kg/GJ)"
$ws.Range("J1").Value = "IPCC_2008"
$ws.Range("I1").Value = "IPCC_1996"
$ws.Range("P1").Value = "IPCC_1996 (kgC/GJ)"

# --- New rows 15-18 (gas names introduce new shared strings next) ---
$ws.Range("C15").Value = "gas_works_gas"
$ws.Range("C16").Value = "coke_oven_gas"
$ws.Range("C17").Value = "blast_furnace_gas"
$ws.Range("C18").Value = "oxygen_steel_furnace_gas"

$ws.Range("M1").Value = "IPCC_96_diff"
$ws.Range("N1").Value = "IPCC_08_diff"

$ws.Range("B15").Value = "1A1b_Mass_Balance_Adjust"
$ws.Range("B16").Value = "1A1b_Mass_Balance_Adjust"
$ws.Range("B17").Value = "1A1b_Mass_Balance_Adjust"
$ws.Range("B18").Value = "1A1b_Mass_Balance_Adjust"

$ws.Range("R1").Value = "kgCO2/GJ"
$ws.Range("Q1").Value = "IPCC_2008"

# --- Remaining (non-shared-string-affecting) cells for rows 15-18 ---
$ws.Range("A15").Value = "default"
$ws.Range("A16").Value = "default"
$ws.Range("A17").Value = "default"
$ws.Range("A18").Value = "default"
$ws.Range("D15").Value = "kt CO2/TJ"
$ws.Range("D16").Value = "kt CO2/TJ"
$ws.Range("D17").Value = "kt CO2/TJ"
$ws.Range("D18").Value = "kt CO2/TJ"

# Row 2
$ws.Range("H2").Formula = '=IF(D2="kt CO2/kJ",E2*1000000*1000000,"")'
$ws.Range("I2").Formula = '=P2*CO2_to_C'
$ws.Range("M2").Formula = '=(H2-I2)/I2'
$ws.Range("P2").Value = 25.8
$ws.Range("Q2").Value = $null
$ws.Range("R2").Value = $null

# Row 3
$ws.Range("H3").Formula = '=IF(D3="kt CO2/kJ",E3*1000000*1000000,"")'
$ws.Range("I3").Formula = '=P3*CO2_to_C'
$ws.Range("K3").Formula = '=S3*CO2_to_C'
$ws.Range("M3").Formula = '=(H3-I3)/I3'
$ws.Range("P3").Formula = '=26.2'
$ws.Range("S3").Value = 27.6

# Row 4
$ws.Range("H4").Formula = '=IF(D4="kt CO2/TJ",E4/1000*1000000,"")'
$ws.Range("I4").Formula = '=P4*CO2_to_C'
$ws.Range("M4").Formula = '=(H4-I4)/I4'
$ws.Range("P4").Formula = '=Q4'
$ws.Range("Q4").Formula = '=R4/CO2_to_C'
$ws.Range("R4").Value = 56.1

# Row 5
$ws.Range("H5").Formula = '=IF(D5="kt CO2/kJ",E5*1000000*1000000,"")'

# Row 6
$ws.Range("H6").Formula = '=IF(D6="kt CO2/TJ",E6/1000*1000000,"")'

# Row 7
$ws.Range("H7").Formula = '=IF(D7="kt CO2/kJ",E7*1000000*1000000,"")'

# Row 8
$ws.Range("H8").Formula = '=IF(D8="kt CO2/kJ",E8*1000000*1000000,"")'

# Row 9
$ws.Range("H9").Formula = '=IF(D9="kt CO2/kJ",E9*1000000*1000000,"")'

# Row 10
$ws.Range("H10").Formula = '=IF(D10="kt CO2/kJ",E10*1000000*1000000,"")'

# Row 11
$ws.Range("H11").Formula = '=IF(D11="kt CO2/kJ",E11*1000000*1000000,"")'

# Row 12
$ws.Range("H12").Formula = '=IF(D12="kt CO2/kJ",E12*1000000*1000000,"")'

# Row 13
$ws.Range("H13").Formula = '=IF(D13="kt CO2/kJ",E13*1000000*1000000,"")'

# Row 14
$ws.Range("H14").Formula = '=IF(D14="kt CO2/kJ",E14*1000000*1000000,"")'

# Row 15 formulas
$ws.Range("E15").Formula = '=Q15*CO2_to_C/1000'
$ws.Range("Q15").Formula = '=R15/CO2_to_C'
$ws.Range("R15").Value = 44.4

# Row 16 formulas
$ws.Range("E16").Formula = '=Q16*CO2_to_C/1000'
$ws.Range("Q16").Formula = '=R16/CO2_to_C'
$ws.Range("R16").Value = 44.4

# Row 17 formulas
$ws.Range("E17").Formula = '=Q17*CO2_to_C/1000'
$ws.Range("Q17").Formula = '=R17/CO2_to_C'
$ws.Range("R17").Value = 260

# Row 18 formulas
$ws.Range("E18").Formula = '=Q18*CO2_to_C/1000'
$ws.Range("Q18").Formula = '=R18/CO2_to_C'
$ws.Range("R18").Value = 182

# --- Number formats / alignment for the new analysis columns ---
# Style 26-equivalent: "0.0" number format, centered
$style26Ranges = @("H2:H14", "I2:I4", "K3", "L3", "O3", "P3:P4", "Q4", "Q15:Q18", "S3")
foreach ($rngAddr in $style26Ranges) {
    $r = $ws.Range($rngAddr)
    $r.NumberFormat = "0.0"
    $r.HorizontalAlignment = -4108
}

# Style 27-equivalent: General number format, centered
$style27Ranges = @("P2", "Q2", "R2", "Q3", "R4", "R15:R18")
foreach ($rngAddr in $style27Ranges) {
    $r = $ws.Range($rngAddr)
    $r.HorizontalAlignment = -4108
}

# Style 28-equivalent: Percent number format, centered
$style28Ranges = @("M2:M4", "N3")
foreach ($rngAddr in $style28Ranges) {
    $r = $ws.Range($rngAddr)
    $r.Style = "Percent"
    $r.HorizontalAlignment = -4108
}

# Style 29-equivalent: 10pt font (header cells)
$ws.Range("M1").Font.Size = 10
$ws.Range("N1").Font.Size = 10

Write-Host "done"
